# Update Active_Outages.xlsx - 6/18/2025, 4:30:02 PM
#
# R1 (sheet1):
#  - G2 "Elapsed Duration(Hrs)" value changes 3929:44:03 -> 3929:44:15
#  - G3 "Elapsed Duration(Hrs)" value changes 69:16:41 -> 69:16:53
#  - New row 6 appended for a new outage entry:
#       Region=R4, Hub Site=LTH0330, Power Source=SCECO+STB,
#       Battery Backup Status=Good, Site Owner=Latis
#       (PCM, Count sites, Fault Level, PCM Created At, Elapsed Duration,
#        Creat Fault First Time Occured and EM Field Feedback are left blank,
#        matching the pattern already used on rows 4 and 5 of this sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R1")

$ws.Range("G2").Value = "3929:44:15"
$ws.Range("G3").Value = "69:16:53"

$ws.Range("B6").Value = "R4"
$ws.Range("D6").Value = "LTH0330"
$ws.Range("I6").Value = "SCECO+STB"
$ws.Range("J6").Value = "Good"
$ws.Range("L6").Value = "Latis"
